# Append a freshly scraped Lancers listing to the "ランサーズ" sheet.
# This mirrors a scraper run at 2026-01-31 01:47:20 that:
#   1. Prepends one new job listing just before the trailing block of rows
#      that were already present (so it lands at row 10, pushing the
#      previous rows 10-12 down to rows 11-13).
#   2. Refreshes the "取得日時" (fetched-at) timestamp on every row, since
#      each scrape run re-stamps the whole sheet.
#   3. Re-creates the hyperlinks on column F so each one still points at
#      the URL that is actually displayed in its row after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-31 01:47:20"

# --- 1. Insert the new row at position 10 (rows 10-12 shift to 11-13) ---
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = $newTimestamp
$ws.Range("B10").Value = "Power Automate(またはGAS)での予約サイト連携フロー構築"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5482835"
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = "◇サイト"

# --- 2. Refresh the timestamp column for every data row (2-13) ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 3. Rebuild the column F hyperlinks so they line up with the rows ---
$ws.Range("F10").Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5460562"
    3  = "https://www.lancers.jp/work/detail/5482607"
    4  = "https://www.lancers.jp/work/detail/5460563"
    5  = "https://www.lancers.jp/work/detail/5475245"
    6  = "https://www.lancers.jp/work/detail/5482462"
    7  = "https://www.lancers.jp/work/detail/5482097"
    8  = "https://www.lancers.jp/work/detail/5482389"
    9  = "https://www.lancers.jp/work/detail/5481859"
    10 = "https://www.lancers.jp/work/detail/5482835"
    11 = "https://www.lancers.jp/work/detail/5418064"
    12 = "https://www.lancers.jp/work/detail/5481715"
    13 = "https://www.lancers.jp/work/detail/5481888"
}

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = $urls[$r]
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Style = "Hyperlink"
}
